$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $origStyle = $c.Style
    $c.Value = "'" + $val
    $c.Style = $origStyle
}

Set-TextValue "D2" "60.589.45"
Set-TextValue "E2" "  -2.34%  "
Set-TextValue "D3" "2.907.26"
Set-TextValue "E3" "  -3.03%  "
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "529.13"
Set-TextValue "E5" "  -4.00%  "
Set-TextValue "E6" "  -7.57%  "
Set-TextValue "E7" "  -0.05%  "
Set-TextValue "D8" "0.555"
Set-TextValue "E8" "  -1.92%  "
Set-TextValue "D9" "2.910.48"
Set-TextValue "E9" "  -2.98%  "
Set-TextValue "E10" "  -3.52%  "
Set-TextValue "D11" "5.93"
Set-TextValue "E11" "  -5.20%  "
Set-TextValue "D12" "0.361"
Set-TextValue "E12" "  -1.19%  "
Set-TextValue "D13" "3.413.14"
Set-TextValue "E13" "  -3.11%  "
Set-TextValue "E14" "  +1.46%  "
Set-TextValue "D15" "60.598.30"
Set-TextValue "E15" "  -2.39%  "
Set-TextValue "D16" "22.65"
Set-TextValue "E16" "  -4.59%  "
Set-TextValue "D17" "2.906.94"
Set-TextValue "E17" "  -3.12%  "
Set-TextValue "E18" "  -3.89%  "
Set-TextValue "D19" "5.05"
Set-TextValue "E19" "  -1.32%  "
Set-TextValue "D20" "11.68"
Set-TextValue "E20" "  -2.46%  "
Set-TextValue "D21" "364.44"
Set-TextValue "E21" "  -7.04%  "
Set-TextValue "D22" "6.62"
Set-TextValue "E22" "  -0.64%  "
Set-TextValue "D23" "0.999"
Set-TextValue "E23" "  -0.13%  "
Set-TextValue "D24" "63.79"
Set-TextValue "E24" "  -2.12%  "
Set-TextValue "D25" "3.022.46"
Set-TextValue "E25" "  -3.31%  "
Set-TextValue "D26" "0.453"
Set-TextValue "E26" "  -3.24%  "
Set-TextValue "E27" "  -2.54%  "
Set-TextValue "E28" "  +0.00%  "
Set-TextValue "E29" "  -7.32%  "
Set-TextValue "E30" "  -8.68%  "
Set-TextValue "E31" "  +0.01%  "
Set-TextValue "D32" "1.68"
Set-TextValue "E32" "  -2.88%  "
Set-TextValue "D33" "19.60"
Set-TextValue "E33" "  -4.33%  "
Set-TextValue "D34" "148.60"
Set-TextValue "E34" "  -6.38%  "
Set-TextValue "D35" "4.37"
Set-TextValue "E35" "  -6.14%  "
Set-TextValue "D36" "5.60"
Set-TextValue "E36" "  -6.97%  "
Set-TextValue "D37" "1.01"
Set-TextValue "E37" "  -6.73%  "
Set-TextValue "E38" "  -6.42%  "
Set-TextValue "D39" "37.96"
Set-TextValue "E39" "  +1.97%  "
Set-TextValue "E40" "  -3.90%  "
Set-TextValue "D41" "2.336.22"
Set-TextValue "E41" "  -4.70%  "
Set-TextValue "D42" "3.68"
Set-TextValue "E42" "  -5.37%  "
Set-TextValue "E43" "  -2.42%  "
Set-TextValue "D44" "20.78"
Set-TextValue "E44" "  -6.95%  "
Set-TextValue "D45" "0.0575"
Set-TextValue "E45" "  -3.18%  "
Set-TextValue "D46" "0.996"
Set-TextValue "E46" "  -0.12%  "
Set-TextValue "D47" "4.97"
Set-TextValue "E47" "  -0.32%  "
Set-TextValue "E48" "  -4.19%  "
Set-TextValue "D49" "0.0934"
Set-TextValue "E49" "  -1.76%  "
Set-TextValue "D50" "10.33"
Set-TextValue "E50" "  -1.16%  "
Set-TextValue "D51" "251.96"
Set-TextValue "E51" "  -5.34%  "
